# coverletterGraphCore.docx edits
#  1. Split the "RE: Application for Internship in Software Frameworks Team"
#     run into four runs (same visible text, new run boundaries).
#  2. Append several new sentences to the end of the "Additionally, I have
#     worked on projects using both C and Java. " paragraph.
#
# Helper: force a run split at absolute document character position $pos
# (i.e. make sure a run boundary exists immediately before the character
# currently sitting at $pos) without altering the character formatting -
# toggling Bold on then back off over [$pos, $endOfEditedSpan) is enough
# to make the host engine keep the runs distinct instead of coalescing
# them back together when it serializes the package.
function Split-RunAt($doc, $pos, $spanEnd) {
    $rr = $doc.Range($pos, $spanEnd)
    $rr.Font.Bold = 1
    $rr.Font.Bold = 0
}

# Insert $parts (an ordered list of strings that concatenate to the text
# actually typed) at the collapsed range $insertionPoint, then recreate
# the run boundaries between every adjacent pair of parts.
function Insert-SplitRuns($doc, $insertionPoint, [string[]]$parts) {
    $full = [string]::Join("", $parts)
    $insertStart = $insertionPoint.Start
    $insertionPoint.Text = $full
    $insertEnd = $insertStart + $full.Length

    $boundaries = New-Object System.Collections.ArrayList
    $offset = $insertStart
    [void]$boundaries.Add($offset)
    for ($i = 0; $i -lt $parts.Length - 1; $i++) {
        $offset += $parts[$i].Length
        [void]$boundaries.Add($offset)
    }

    # Apply from the last boundary back to the first so that each toggle
    # only ever touches the not-yet-split tail of the inserted text.
    for ($i = $boundaries.Count - 1; $i -ge 0; $i--) {
        Split-RunAt $doc $boundaries[$i] $insertEnd
    }
    return $insertEnd
}

$d = $word.ActiveDocument

# --- Edit 1: split the "RE: Application for Internship in Software Frameworks Team" run ---

$rIntern = $d.Content.Duplicate
[void]$rIntern.Find.Execute("Intern", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rIntern.Font.Bold = 1
$rIntern.Font.Bold = 0

$rShip = $d.Content.Duplicate
[void]$rShip.Find.Execute("ship in ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rShip.Font.Bold = 1
$rShip.Font.Bold = 0

# --- Edit 2: append the new sentences about C/C++/Java experience ---

$anchor = $d.Content.Duplicate
[void]$anchor.Find.Execute("worked on projects using both C and Java. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)

$newParts = @(
    "I",
    "mplementing linked lists and ",
    "a basic sketching tool in C, as well as several object-oriented applications including a board game and multithreaded 3D model analyzer in Java, ",
    "so the syntax and object-orientated paradigm ",
    "of C++ ",
    "will be very familiar to me",
    ". ",
    "I have taught myself the electric guitar in my",
    " spare time,",
    " so ",
    "I am no stranger to learning new skills in my own time and would pick up C++ quickly."
)

[void](Insert-SplitRuns $d $anchor $newParts)

Write-Output "edits applied"
